$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasting Notes")

# Add new tasting-note row (row 21): date, shared day-# formula, score, notes.
$ws.Range("A21").Value = 44212
$ws.Range("B21").Formula = "=A21-`$A`$6"
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = "cooled and served @ 11 C, light foam which dissolves quickly, very high carbonation, fresh dry mouthfeel and sweet/sour taste."

# Match formatting of the row above (row 20) for the new row.
$ws.Range("A20:D20").Copy()
$ws.Range("A21:D21").PasteSpecial(-4122)
$ws.Range("A21").Value = 44212
$ws.Range("B21").Formula = "=A21-`$A`$6"
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = "cooled and served @ 11 C, light foam which dissolves quickly, very high carbonation, fresh dry mouthfeel and sweet/sour taste."

$ws.Application.CutCopyMode = 0

# Update the view so the new row is reachable, matching the saved workbook state.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("D22").Select()
